$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.023.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.111.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.14%  "
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("E13").Value = "  +6.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.00%  "
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.635.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.34%  "
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.114.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "62.911.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "456.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.24%  "
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("E23").Value = "  +6.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +5.46%  "
$ws.Range("E28").Value = "  +7.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.40%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.99%  "
$ws.Range("E32").Value = "  +16.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.40%  "
$ws.Range("E35").Value = "  +6.92%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("E37").Value = "  +6.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "428.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.941.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0373"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.33%  "
